$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 153, pushing existing rows 153:184 down to 154:185
$ws.Rows.Item(153).Insert()

# Populate the newly inserted row 153 with the new weekly record
$ws.Cells.Item(153, 1).Value = 5
$ws.Cells.Item(153, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(153, 3).Value = "Maule"
$ws.Cells.Item(153, 4).Value = 44551
$ws.Cells.Item(153, 5).Value = 7
$ws.Cells.Item(153, 6).Value = 100112008
$ws.Cells.Item(153, 7).Value = "Coliflor"
$ws.Cells.Item(153, 8).Value = "Sin especificar"
$ws.Cells.Item(153, 9).Value = "Primera"
$ws.Cells.Item(153, 10).Value = 3000
$ws.Cells.Item(153, 11).Value = 900
$ws.Cells.Item(153, 12).Value = 900
$ws.Cells.Item(153, 13).Value = 900
$ws.Cells.Item(153, 14).Value = "`$/unidad"
$ws.Cells.Item(153, 15).Value = "Región del Maule"
$ws.Cells.Item(153, 16).Value = 900
$ws.Cells.Item(153, 17).Value = 1
$ws.Cells.Item(153, 18).Value = "Hortaliza"
